# "took out the summary page" — removes the per-row "View Job" hyperlink
# cells from the Jobs sheet (keeping only the last row, whose text becomes
# the literal URL) and strips the link styling/relationships from the
# Summary sheet's Top Jobs table (again keeping only the last row's link),
# plus drops the now-unused "Category Filter" row from the Summary sheet.
#
# NOTE: in this COM host, Range.Hyperlinks.Delete() removes every
# hyperlink on the parent worksheet (not just the ones in the range), so
# each sheet's links are cleared once up front and the single link that
# should survive is re-added afterwards with Hyperlinks.Add.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Jobs": column L ("Link") rows 2-19 currently hold a styled
# "View Job" hyperlink cell. Remove the hyperlink + cell entirely for
# rows 2-18; row 19 keeps a hyperlink but the visible text becomes the
# raw URL instead of "View Job".
# ---------------------------------------------------------------------
$jobs = $wb.Worksheets.Item("Jobs")

$jobs.Hyperlinks.Delete()

for ($row = 2; $row -le 18; $row++) {
    $jobs.Cells.Item($row, 12).Clear()
}

$jobsLastUrl = "https://remotive.com/remote-jobs/marketing/senior-amazon-brand-manager-2082736"
$jobs.Range("L19").Value = $jobsLastUrl
$jobs.Hyperlinks.Add($jobs.Range("L19"), $jobsLastUrl)
$jobs.Range("L19").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Sheet "Summary": drop the "Category Filter" row (row 6) completely.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A6:B6").Clear()

# ---------------------------------------------------------------------
# Sheet "Summary": Top Jobs table, column E ("Link") rows 14-23. Strip
# the hyperlink + Hyperlink style for rows 14-22 and replace the
# "View Job" text with the raw URL; row 23 keeps a hyperlink/style but
# also swaps its text to the raw URL.
# ---------------------------------------------------------------------
$summary.Hyperlinks.Delete()

$summaryLinks = @{
    14 = "https://remotive.com/remote-jobs/project-management/order-management-and-operations-manager-2088635"
    15 = "https://remotive.com/remote-jobs/devops/senior-devops-engineer-2070150"
    16 = "https://remotive.com/remote-jobs/software-development/full-stack-developer-6-months-extendable-2088631"
    17 = "https://remotive.com/remote-jobs/ai-ml/ai-native-cloud-infrastructure-generalist-m-f-d-2088634"
    18 = "https://remotive.com/remote-jobs/software-development/tech-lead-databricks-data-engineer-2069747"
    19 = "https://remotive.com/remote-jobs/software-development/senior-python-backend-developer-2088624"
    20 = "https://remotive.com/remote-jobs/software-development/senior-independent-ai-engineer-architect-1919266"
    21 = "https://remotive.com/remote-jobs/software-development/senior-independent-software-developer-1919265"
    22 = "https://remotive.com/remote-jobs/customer-service/client-support-specialist-2086826"
}

foreach ($row in $summaryLinks.Keys) {
    $cell = $summary.Cells.Item($row, 5)
    $cell.ClearFormats()
    $cell.Value = $summaryLinks[$row]
}

$summaryLastUrl = "https://remotive.com/remote-jobs/software-development/tech-lead-full-stack-rails-engineer-2069746"
$summary.Range("E23").Value = $summaryLastUrl
$summary.Hyperlinks.Add($summary.Range("E23"), $summaryLastUrl)
$summary.Range("E23").Style = "Hyperlink"
